$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New yearly columns: S = 2021, T = 2022, added for every data row (4-14),
# mirroring the formatting already used by column R (the previous last year).
$values = @{
    4  = @(2021, 2022)
    5  = @(2.5, 2.6)
    6  = @(2.5, 1.8)
    7  = @(1.6, 2.6)
    8  = @(3.6, 1.9)
    9  = @(5.8, 3.9)
    10 = @(1.1, 3.2)
    11 = @(1.1, 3.3)
    12 = @(5.1, 2.5)
    13 = @(2.3, 1.9)
    14 = @(2.1, 2.5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]

    # Copy column R's cell format (style) into S and T before writing values,
    # so the new cells carry the same number format / borders as the rest of
    # the table row.
    $null = $ws.Cells.Item($row, 18).Copy()
    $null = $ws.Cells.Item($row, 19).PasteSpecial(-4122)
    $null = $ws.Cells.Item($row, 20).PasteSpecial(-4122)

    $ws.Cells.Item($row, 19).Value = $pair[0]
    $ws.Cells.Item($row, 20).Value = $pair[1]
}

# Match the saved selection state recorded in the diff.
$null = $ws.Range("V7").Select()
